$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.832.14"
$ws.Range("E2").Value = "  -2.76%  "

$ws.Range("D3").Value = "3.179.30"
$ws.Range("E3").Value = "  -1.65%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.10"
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.83"
$ws.Range("E6").Value = "  -4.13%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.178.34"
$ws.Range("E8").Value = "  -1.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  -3.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  -4.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.59"
$ws.Range("E11").Value = "  -1.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.476"
$ws.Range("E12").Value = "  -6.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("E13").Value = "  -4.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.03"
$ws.Range("E14").Value = "  -5.37%  "

$ws.Range("D15").Value = "3.699.33"
$ws.Range("E15").Value = "  -1.70%  "

$ws.Range("D16").Value = "64.850.21"
$ws.Range("E16").Value = "  -2.77%  "

$ws.Range("D17").Value = "3.178.33"
$ws.Range("E17").Value = "  -1.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.114"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.04"
$ws.Range("E19").Value = "  -4.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.41"
$ws.Range("E20").Value = "  -5.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.84"
$ws.Range("E21").Value = "  -2.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.717"
$ws.Range("E22").Value = "  -2.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.78"
$ws.Range("E23").Value = "  -3.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.96"
$ws.Range("E24").Value = "  -4.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.49"
$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.93"
$ws.Range("E27").Value = "  -2.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.73"
$ws.Range("E28").Value = "  -4.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  -5.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.09"
$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("E31").Value = "  +6.51%  "

$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.73"
$ws.Range("E32").Value = "  -8.01%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.93"
$ws.Range("E34").Value = "  -4.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.12"
$ws.Range("E35").Value = "  -5.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.16"
$ws.Range("E36").Value = "  -5.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.76"
$ws.Range("E37").Value = "  -1.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  +5.91%  "

$ws.Range("D39").Value = "0.0₃0743"
$ws.Range("E39").Value = "  -4.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "464.25"
$ws.Range("E40").Value = "  -8.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.126"
$ws.Range("E41").Value = "  -2.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0404"
$ws.Range("E42").Value = "  -4.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.52"
$ws.Range("E43").Value = "  -2.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.45"
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").Value = "2.898.28"
$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.276"
$ws.Range("E46").Value = "  -7.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.43"
$ws.Range("E47").Value = "  -3.54%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.35"
$ws.Range("E50").Value = "  -2.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.57"
$ws.Range("E51").Value = "  -1.82%  "
